# Updates the "cryptos" price/volume table to the latest scraped values.
# Note: several "Price" (column D) values look numeric (e.g. "354.32") but
# must remain stored as text, exactly like the rest of the sheet. Setting
# them through .Value with a leading apostrophe forces text entry (as the
# Excel UI would), and resetting .Style back to "Normal" afterwards avoids
# leaving behind an extra quote-prefix cell style that isn't part of the
# intended edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.051.40'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.978.53'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'354.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'112.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").Value = "'0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = "'0.633"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("D11").Value = "'0.0898"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.96%  '
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").Value = "'20.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '3.449.58'
$ws.Range("E15").Value = '  +3.89%  '
$ws.Range("D16").Value = '2.986.51'
$ws.Range("E16").Value = '  +2.39%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '52.192.12'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = "'7.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = "'3.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").Value = "'14.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.10%  '
$ws.Range("D22").Value = '0.0₃0992'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = "'71.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = "'270.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("E26").Value = '  +9.82%  '
$ws.Range("D27").Value = "'27.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = "'7.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +22.05%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = "'0.116"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +28.88%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = "'10.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").Value = "'37.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.56%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = "'6.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.28%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = "'53.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = "'0.0452"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").Value = "'2.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.14%  '
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = "'3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.23%  '
$ws.Range("D39").Value = "'19.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = "'2.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("D41").Value = "'2.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.50%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = "'23.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.67%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = "'0.118"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = "'2.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D45").Value = "'3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").Value = '2.179.13'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = "'114.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.85%  '
$ws.Range("D49").Value = "'0.245"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = "'0.0345"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.24%  '
$ws.Range("D51").Value = "'0.945"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.17%  '
